$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(15, 8).Value = 696357.8
$ws.Cells.Item(15, 9).Value = 696357.8
$ws.Cells.Item(15, 11).Value = 2089073.4
$ws.Cells.Item(15, 13).Value = -2088904.4
$ws.Cells.Item(17, 8).Value = 2357.1133
$ws.Cells.Item(17, 10).Value = 2357.1133
$ws.Cells.Item(17, 12).Value = 7071.3399
$ws.Cells.Item(17, 14).Value = -7407.3399
$ws.Cells.Item(18, 8).Value = 1671.9333
$ws.Cells.Item(18, 9).Value = 356.75
$ws.Cells.Item(18, 11).Value = 356.75
$ws.Cells.Item(18, 13).Value = -72.75
$ws.Cells.Item(19, 8).Value = 936.41174
$ws.Cells.Item(19, 10).Value = 1035
$ws.Cells.Item(19, 12).Value = 1035
$ws.Cells.Item(19, 14).Value = -1385
$ws.Cells.Item(48, 8).Value = 7399.6
$ws.Cells.Item(48, 10).Value = 7166
$ws.Cells.Item(48, 12).Value = 21498
$ws.Cells.Item(48, 14).Value = -22082
$ws.Cells.Item(56, 8).Value = 7399.6
$ws.Cells.Item(56, 10).Value = 7166
$ws.Cells.Item(56, 12).Value = 21498
$ws.Cells.Item(56, 14).Value = -22566
$ws.Cells.Item(86, 8).Value = 6140.7
$ws.Cells.Item(86, 9).Value = 5866.32
$ws.Cells.Item(86, 10).Value = 6598
$ws.Cells.Item(86, 11).Value = 5866.32
$ws.Cells.Item(86, 12).Value = 6598
$ws.Cells.Item(86, 13).Value = -4743.32
$ws.Cells.Item(86, 14).Value = -8844
$ws.Cells.Item(89, 8).Value = 6140.7
$ws.Cells.Item(89, 9).Value = 5866.32
$ws.Cells.Item(89, 10).Value = 6598
$ws.Cells.Item(89, 11).Value = 29331.6
$ws.Cells.Item(89, 12).Value = 32990
$ws.Cells.Item(89, 13).Value = -23715.6
$ws.Cells.Item(89, 14).Value = -44222
$ws.Cells.Item(103, 8).Value = 298
$ws.Cells.Item(103, 9).Value = 122.5
$ws.Cells.Item(103, 10).Value = 1000
$ws.Cells.Item(103, 11).Value = 367.5
$ws.Cells.Item(103, 12).Value = 3000
$ws.Cells.Item(103, 13).Value = 218.5
$ws.Cells.Item(103, 14).Value = -4172
$ws.Cells.Item(111, 8).Value = 71136.31
$ws.Cells.Item(111, 9).Value = 698.6667
$ws.Cells.Item(111, 11).Value = 2096.0001
$ws.Cells.Item(111, 13).Value = 970.9998999999998
$ws.Cells.Item(112, 8).Value = 4780
$ws.Cells.Item(112, 10).Value = 5270
$ws.Cells.Item(112, 12).Value = 15810
$ws.Cells.Item(112, 14).Value = -18026
$ws.Cells.Item(113, 8).Value = 4689.52
$ws.Cells.Item(113, 9).Value = 3637.7058
$ws.Cells.Item(113, 10).Value = 6924.625
$ws.Cells.Item(113, 11).Value = 3637.7058
$ws.Cells.Item(113, 12).Value = 6924.625
$ws.Cells.Item(113, 13).Value = -383.7058000000002
$ws.Cells.Item(113, 14).Value = -13432.625
$ws.Cells.Item(127, 8).Value = 5189
$ws.Cells.Item(127, 9).Value = 5189
$ws.Cells.Item(127, 11).Value = 15567
$ws.Cells.Item(127, 13).Value = -10607
$ws.Cells.Item(137, 8).Value = 1922.65
$ws.Cells.Item(137, 9).Value = 1729.3611
$ws.Cells.Item(137, 10).Value = 3662.25
$ws.Cells.Item(137, 11).Value = 5188.0833
$ws.Cells.Item(137, 12).Value = 10986.75
$ws.Cells.Item(137, 13).Value = -2638.0833
$ws.Cells.Item(137, 14).Value = -16086.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 10053.583
$ws.Cells.Item(61, 9).Value = 10735
$ws.Cells.Item(61, 11).Value = 10735
$ws.Cells.Item(61, 13).Value = -10523
$ws.Cells.Item(74, 8).Value = 5970.15
$ws.Cells.Item(74, 9).Value = 2086.5
$ws.Cells.Item(74, 11).Value = 2086.5
$ws.Cells.Item(74, 13).Value = -1212.5
$ws.Cells.Item(77, 8).Value = 5970.15
$ws.Cells.Item(77, 9).Value = 2086.5
$ws.Cells.Item(77, 11).Value = 10432.5
$ws.Cells.Item(77, 13).Value = -6064.5
$ws.Cells.Item(102, 8).Value = 1606.2222
$ws.Cells.Item(102, 9).Value = 1443.25
$ws.Cells.Item(102, 10).Value = 2910
$ws.Cells.Item(102, 11).Value = 1443.25
$ws.Cells.Item(102, 12).Value = 2910
$ws.Cells.Item(102, 13).Value = 178.75
$ws.Cells.Item(102, 14).Value = -6154
$ws.Cells.Item(110, 8).Value = 1980.5
$ws.Cells.Item(110, 9).Value = 1965.2554
$ws.Cells.Item(110, 11).Value = 1965.2554
$ws.Cells.Item(110, 13).Value = 79.74459999999999
$ws.Cells.Item(122, 8).Value = 3182.75
$ws.Cells.Item(122, 9).Value = 1438.1666
$ws.Cells.Item(122, 11).Value = 4314.4998
$ws.Cells.Item(122, 13).Value = -1864.4998
$ws.Cells.Item(132, 8).Value = 4061.2856
$ws.Cells.Item(132, 9).Value = 2007.1923
$ws.Cells.Item(132, 11).Value = 6021.5769
$ws.Cells.Item(132, 13).Value = -3491.5769
$ws.Cells.Item(135, 8).Value = 44775.75
$ws.Cells.Item(135, 10).Value = 44775.75
$ws.Cells.Item(135, 12).Value = 44775.75
$ws.Cells.Item(135, 14).Value = -54915.75
$ws.Cells.Item(136, 8).Value = 10053.583
$ws.Cells.Item(136, 9).Value = 10735
$ws.Cells.Item(136, 11).Value = 32205
$ws.Cells.Item(136, 13).Value = -29655

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(94, 8).Value = 588.26086
$ws.Cells.Item(94, 9).Value = 591.95
$ws.Cells.Item(94, 10).Value = 563.6667
$ws.Cells.Item(94, 11).Value = 591.95
$ws.Cells.Item(94, 12).Value = 563.6667
$ws.Cells.Item(94, 13).Value = -140.95
$ws.Cells.Item(94, 14).Value = -1465.6667
$ws.Cells.Item(99, 8).Value = 2184.8333
$ws.Cells.Item(99, 9).Value = 2021.9
$ws.Cells.Item(99, 11).Value = 2021.9
$ws.Cells.Item(99, 13).Value = -523.9000000000001
$ws.Cells.Item(107, 8).Value = 4008.8
$ws.Cells.Item(107, 9).Value = 3974.75
$ws.Cells.Item(107, 11).Value = 3974.75
$ws.Cells.Item(107, 13).Value = -2054.75
$ws.Cells.Item(134, 8).Value = 4016.9167
$ws.Cells.Item(134, 9).Value = 2223.524
$ws.Cells.Item(134, 11).Value = 6670.572
$ws.Cells.Item(134, 13).Value = -4135.572
$ws.Cells.Item(138, 8).Value = 95000
$ws.Cells.Item(138, 10).Value = 95000
$ws.Cells.Item(138, 12).Value = 95000
$ws.Cells.Item(138, 14).Value = -105280

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 861
$ws.Cells.Item(16, 9).Value = 432.57144
$ws.Cells.Item(16, 11).Value = 432.57144
$ws.Cells.Item(16, 13).Value = -145.57144
$ws.Cells.Item(31, 8).Value = 8123.9644
$ws.Cells.Item(31, 9).Value = 3415.111
$ws.Cells.Item(31, 11).Value = 3415.111
$ws.Cells.Item(31, 13).Value = -3120.111
$ws.Cells.Item(34, 8).Value = 8123.9644
$ws.Cells.Item(34, 9).Value = 3415.111
$ws.Cells.Item(34, 11).Value = 3415.111
$ws.Cells.Item(34, 13).Value = -3213.111
$ws.Cells.Item(58, 8).Value = 3334.0286
$ws.Cells.Item(58, 9).Value = 2498.1428
$ws.Cells.Item(58, 10).Value = 4587.857
$ws.Cells.Item(58, 11).Value = 2498.1428
$ws.Cells.Item(58, 12).Value = 4587.857
$ws.Cells.Item(58, 13).Value = -2295.1428
$ws.Cells.Item(58, 14).Value = -4993.857
$ws.Cells.Item(62, 8).Value = 3712.9412
$ws.Cells.Item(62, 9).Value = 3348.4285
$ws.Cells.Item(62, 11).Value = 3348.4285
$ws.Cells.Item(62, 13).Value = -2724.4285
$ws.Cells.Item(65, 8).Value = 3712.9412
$ws.Cells.Item(65, 9).Value = 3348.4285
$ws.Cells.Item(65, 11).Value = 16742.1425
$ws.Cells.Item(65, 13).Value = -13622.1425
$ws.Cells.Item(74, 8).Value = 57666.445
$ws.Cells.Item(74, 10).Value = 72999.336
$ws.Cells.Item(74, 12).Value = 72999.336
$ws.Cells.Item(74, 14).Value = -74747.336
$ws.Cells.Item(77, 8).Value = 57666.445
$ws.Cells.Item(77, 10).Value = 72999.336
$ws.Cells.Item(77, 12).Value = 218998.008
$ws.Cells.Item(77, 14).Value = -227734.008
$ws.Cells.Item(99, 8).Value = 2431.926
$ws.Cells.Item(99, 9).Value = 2463.182
$ws.Cells.Item(99, 10).Value = 2294.4
$ws.Cells.Item(99, 11).Value = 2463.182
$ws.Cells.Item(99, 12).Value = 2294.4
$ws.Cells.Item(99, 13).Value = -965.1819999999998
$ws.Cells.Item(99, 14).Value = -5290.4
$ws.Cells.Item(105, 8).Value = 3773.1
$ws.Cells.Item(105, 9).Value = 3378.875
$ws.Cells.Item(105, 10).Value = 5350
$ws.Cells.Item(105, 11).Value = 3378.875
$ws.Cells.Item(105, 12).Value = 5350
$ws.Cells.Item(105, 13).Value = -1631.875
$ws.Cells.Item(105, 14).Value = -8844
$ws.Cells.Item(107, 8).Value = 1075.1666
$ws.Cells.Item(107, 9).Value = 968.5217
$ws.Cells.Item(107, 10).Value = 1425.5714
$ws.Cells.Item(107, 11).Value = 968.5217
$ws.Cells.Item(107, 12).Value = 1425.5714
$ws.Cells.Item(107, 13).Value = 951.4783
$ws.Cells.Item(107, 14).Value = -5265.5714
$ws.Cells.Item(113, 8).Value = 861
$ws.Cells.Item(113, 9).Value = 432.57144
$ws.Cells.Item(113, 11).Value = 432.57144
$ws.Cells.Item(113, 13).Value = 1737.42856
$ws.Cells.Item(126, 8).Value = 2431.926
$ws.Cells.Item(126, 9).Value = 2463.182
$ws.Cells.Item(126, 10).Value = 2294.4
$ws.Cells.Item(126, 11).Value = 7389.545999999999
$ws.Cells.Item(126, 12).Value = 6883.200000000001
$ws.Cells.Item(126, 13).Value = -4919.545999999999
$ws.Cells.Item(126, 14).Value = -11823.2
$ws.Cells.Item(134, 8).Value = 5222.5
$ws.Cells.Item(134, 9).Value = 4312.2383
$ws.Cells.Item(134, 10).Value = 11594.333
$ws.Cells.Item(134, 11).Value = 12936.7149
$ws.Cells.Item(134, 12).Value = 34782.999
$ws.Cells.Item(134, 13).Value = -10401.7149
$ws.Cells.Item(134, 14).Value = -39852.999
$ws.Cells.Item(136, 8).Value = 3334.0286
$ws.Cells.Item(136, 9).Value = 2498.1428
$ws.Cells.Item(136, 10).Value = 4587.857
$ws.Cells.Item(136, 11).Value = 7494.428400000001
$ws.Cells.Item(136, 12).Value = 13763.571
$ws.Cells.Item(136, 13).Value = -4944.428400000001
$ws.Cells.Item(136, 14).Value = -18863.571

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 49616384
$ws.Cells.Item(4, 10).Value = 1112444.4
$ws.Cells.Item(4, 12).Value = 3337333.2
$ws.Cells.Item(4, 14).Value = -3337557.2
$ws.Cells.Item(56, 8).Value = 7569.7144
$ws.Cells.Item(56, 9).Value = 7569.7144
$ws.Cells.Item(56, 11).Value = 7569.7144
$ws.Cells.Item(56, 13).Value = -7039.7144

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(10, 8).Value = 58983.332
$ws.Cells.Item(10, 9).Value = 20000
$ws.Cells.Item(10, 10).Value = 66780
$ws.Cells.Item(10, 11).Value = 20000
$ws.Cells.Item(10, 12).Value = 66780
$ws.Cells.Item(10, 13).Value = -19831
$ws.Cells.Item(10, 14).Value = -67118
$ws.Cells.Item(41, 9).Value = 10000
$ws.Cells.Item(41, 11).Value = 10000
$ws.Cells.Item(41, 13).Value = -9645
$ws.Cells.Item(70, 8).Value = 6253.6665
$ws.Cells.Item(70, 10).Value = 6291.8335
$ws.Cells.Item(70, 12).Value = 6291.8335
$ws.Cells.Item(70, 14).Value = -6831.8335
$ws.Cells.Item(73, 8).Value = 6253.6665
$ws.Cells.Item(73, 10).Value = 6291.8335
$ws.Cells.Item(73, 12).Value = 6291.8335
$ws.Cells.Item(73, 14).Value = -8163.8335
$ws.Cells.Item(97, 8).Value = 1049.4286
$ws.Cells.Item(97, 9).Value = 337.25
$ws.Cells.Item(97, 10).Value = 1999
$ws.Cells.Item(97, 11).Value = 337.25
$ws.Cells.Item(97, 12).Value = 1999
$ws.Cells.Item(97, 13).Value = 158.75
$ws.Cells.Item(97, 14).Value = -2991
$ws.Cells.Item(109, 8).Value = 62500
$ws.Cells.Item(109, 10).Value = 62500
$ws.Cells.Item(109, 12).Value = 62500
$ws.Cells.Item(109, 14).Value = -64580
$ws.Cells.Item(122, 8).Value = 5202.4443
$ws.Cells.Item(122, 9).Value = 1012.25
$ws.Cells.Item(122, 11).Value = 3036.75
$ws.Cells.Item(122, 13).Value = -586.75
$ws.Cells.Item(126, 8).Value = 5923.077
$ws.Cells.Item(126, 9).Value = 4111.1113
$ws.Cells.Item(126, 10).Value = 10000
$ws.Cells.Item(126, 11).Value = 12333.3339
$ws.Cells.Item(126, 12).Value = 30000
$ws.Cells.Item(126, 13).Value = -9863.333899999998
$ws.Cells.Item(126, 14).Value = -34940

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(61, 8).Value = 3480.3076
$ws.Cells.Item(61, 9).Value = 2973.9
$ws.Cells.Item(61, 10).Value = 5168.3335
$ws.Cells.Item(61, 11).Value = 2973.9
$ws.Cells.Item(61, 12).Value = 5168.3335
$ws.Cells.Item(61, 13).Value = -2771.9
$ws.Cells.Item(61, 14).Value = -5572.3335
$ws.Cells.Item(68, 8).Value = 3999.6
$ws.Cells.Item(68, 9).Value = 3999.5
$ws.Cells.Item(68, 11).Value = 3999.5
$ws.Cells.Item(68, 13).Value = -3250.5
$ws.Cells.Item(71, 8).Value = 3999.6
$ws.Cells.Item(71, 9).Value = 3999.5
$ws.Cells.Item(71, 11).Value = 19997.5
$ws.Cells.Item(71, 13).Value = -16253.5
$ws.Cells.Item(100, 8).Value = 2411.5
$ws.Cells.Item(100, 9).Value = 2293.8
$ws.Cells.Item(100, 10).Value = 3000
$ws.Cells.Item(100, 11).Value = 2293.8
$ws.Cells.Item(100, 12).Value = 3000
$ws.Cells.Item(100, 13).Value = -1752.8
$ws.Cells.Item(100, 14).Value = -4082
$ws.Cells.Item(113, 8).Value = 3480.3076
$ws.Cells.Item(113, 9).Value = 2973.9
$ws.Cells.Item(113, 10).Value = 5168.3335
$ws.Cells.Item(113, 11).Value = 2973.9
$ws.Cells.Item(113, 12).Value = 5168.3335
$ws.Cells.Item(113, 13).Value = -803.9000000000001
$ws.Cells.Item(113, 14).Value = -9508.333500000001
$ws.Cells.Item(122, 8).Value = 5047.3335
$ws.Cells.Item(122, 9).Value = 4178.25
$ws.Cells.Item(122, 11).Value = 12534.75
$ws.Cells.Item(122, 13).Value = -10084.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(96, 8).Value = 6124.9165
$ws.Cells.Item(96, 9).Value = 1763.25
$ws.Cells.Item(96, 10).Value = 8305.75
$ws.Cells.Item(96, 11).Value = 1763.25
$ws.Cells.Item(96, 12).Value = 8305.75
$ws.Cells.Item(96, 13).Value = -390.25
$ws.Cells.Item(96, 14).Value = -11051.75
$ws.Cells.Item(132, 8).Value = 2281.5557
$ws.Cells.Item(132, 9).Value = 1736.5116
$ws.Cells.Item(132, 11).Value = 5209.5348
$ws.Cells.Item(132, 13).Value = -2679.5348
